# Lesson 2/5 homework tweak:
#  "Вариант №N" -> "Шаг №N" (with slightly different run-splitting / bookmark
#  / spell-check markup depending on which of the 6 occurrences it is), and
#  the stray "_GoBack" bookmark that used to sit on the "Добавьте свойство"
#  bullet is removed (it effectively "moved" to sit inside the very first
#  "Шаг №1" heading instead).

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-RunsXml([object]$range, [string]$innerXml) {
    # Replaces the contents of $range (which must NOT include the
    # paragraph mark) with the supplied run-level OOXML, preserving the
    # paragraph's own properties (pPr), rsid attributes, etc.
    $full = $pkgHeader + '<w:p>' + $innerXml + '</w:p>' + $pkgFooter
    $range.InsertXML($full)
}

# Remove the old "_GoBack" bookmark that used to sit on the "Добавьте
# свойство ..." bullet point *before* inserting the new one further up in
# the document (it effectively "moves" into the first "Шаг №1" heading) -
# this has to happen first, otherwise the freshly inserted bookmark would
# itself get wiped out since both share the same reserved name.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Walk the document paragraphs in order and patch each of the 6
# "Вариант №N" heading paragraphs according to its position among them.
$count = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $txt = $p.Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt.StartsWith("Вариант")) {
        $count = $count + 1

        $r = $p.Range
        $rEnd = $r.End - 1   # exclude the paragraph mark
        $rStart = $r.Start

        if ($txt -eq "Вариант №1" -and $count -eq 1) {
            # Occurrence 1: split "Вариант №1" -> "Ша" | bookmark _GoBack | "г" | " №1"
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:r><w:t>Ша</w:t></w:r>' + `
                     '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
                     '<w:r><w:t>г</w:t></w:r>' + `
                     '<w:r><w:t xml:space="preserve"> №1</w:t></w:r>'
            Set-RunsXml $target $inner
        }
        elseif ($txt -eq "Вариант №2" -and $count -eq 2) {
            # Occurrence 2: "Вариант №2" -> "Шаг" | " " | "№2"
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:r><w:t>Шаг</w:t></w:r>' + `
                     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
                     '<w:r><w:t>№2</w:t></w:r>'
            Set-RunsXml $target $inner
        }
        elseif ($txt -eq "Вариант №3" -and $count -eq 3) {
            # Occurrence 3: the "Вариант №" run is replaced by three runs;
            # the trailing "3" run keeps its own rsid and is re-emitted
            # unchanged (re-stating it, instead of leaving the original
            # run in place, avoids the target range/run-splice ordering
            # this host's InsertXML otherwise produces).
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:r><w:t>Шаг</w:t></w:r>' + `
                     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
                     '<w:r><w:t>№</w:t></w:r>' + `
                     '<w:r w:rsidR="004B0C3C"><w:t>3</w:t></w:r>'
            Set-RunsXml $target $inner
        }
        elseif ($txt -eq "Вариант №1" -and $count -eq 4) {
            # Occurrence 4: spell-checked "Шаг" run (lang en-US) + " " + "№1"
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:proofErr w:type="spellStart"/>' + `
                     '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Шаг</w:t></w:r>' + `
                     '<w:proofErr w:type="spellEnd"/>' + `
                     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
                     '<w:r><w:t>№1</w:t></w:r>'
            Set-RunsXml $target $inner
        }
        elseif ($txt -eq "Вариант №2" -and $count -eq 5) {
            # Occurrence 5: spell-checked "Шаг" run (lang en-US) + " " + "№2"
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:proofErr w:type="spellStart"/>' + `
                     '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Шаг</w:t></w:r>' + `
                     '<w:proofErr w:type="spellEnd"/>' + `
                     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
                     '<w:r><w:t>№2</w:t></w:r>'
            Set-RunsXml $target $inner
        }
        elseif ($txt -eq "Вариант №3" -and $count -eq 6) {
            # Occurrence 6: same idea as occurrence 3, but the trailing "3"
            # run here carries w:rsidRPr="00B960C2" instead.
            $target = $d.Range($rStart, $rEnd)
            $inner = '<w:r><w:t>Шаг</w:t></w:r>' + `
                     '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
                     '<w:r><w:t>№</w:t></w:r>' + `
                     '<w:r w:rsidRPr="00B960C2"><w:t>3</w:t></w:r>'
            Set-RunsXml $target $inner
        }
    }
}
